$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data: STAB / 3 / Stability Test Report / TRUE
$ws.Range("A4").Value = "STAB"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "Stability Test Report"
$ws.Range("D4").Value = $true

# Update the selected cell/range to match the final saved view state
$ws.Range("D8").Select()
